# Update template's ontology terms
# - Fill in unit / term-source / term-accession columns (C, D, F, G, H, J, K, M, N, P, Q, R, T, U)
#   for the "molecule" rows of the GoodDingo21 annotation table on sheet "2EXT01_RNA"
# - Change "Protein" (row 5, molecule list) to lower-case "protein" and give it its real
#   NFDI4PSO ontology term source/accession
# - Bump template version 1.1.4 -> 1.1.5 on SwateTemplateMetadata

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2EXT01_RNA")

# --- Row 2 (RNA (Transcriptomics)) ---
$ws.Range("C2").Value = "user-specific"
$ws.Range("D2").Value = "user-specific"
$ws.Range("F2").Value = "milligram"
$ws.Range("G2").Value = "UO"
$ws.Range("H2").Value = "http://purl.obolibrary.org/obo/UO_0000022"
$ws.Range("J2").Value = "user-specific"
$ws.Range("K2").Value = "user-specific"
$ws.Range("M2").Value = "user-specific"
$ws.Range("N2").Value = "user-specific"
$ws.Range("P2").Value = "microliter"
$ws.Range("Q2").Value = "UO"
$ws.Range("R2").Value = "http://purl.obolibrary.org/obo/UO_0000101"
$ws.Range("T2").Value = "user-specific"
$ws.Range("U2").Value = "user-specific"

# --- Row 3 (DNA (Genomics)) ---
$ws.Range("C3").Value = "user-specific"
$ws.Range("D3").Value = "user-specific"
$ws.Range("F3").Value = "milligram"
$ws.Range("G3").Value = "UO"
$ws.Range("H3").Value = "http://purl.obolibrary.org/obo/UO_0000022"
$ws.Range("J3").Value = "user-specific"
$ws.Range("K3").Value = "user-specific"
$ws.Range("M3").Value = "user-specific"
$ws.Range("N3").Value = "user-specific"
$ws.Range("P3").Value = "microliter"
$ws.Range("Q3").Value = "UO"
$ws.Range("R3").Value = "http://purl.obolibrary.org/obo/UO_0000101"
$ws.Range("T3").Value = "user-specific"
$ws.Range("U3").Value = "user-specific"

# --- Row 4 (Metabolites) ---
$ws.Range("C4").Value = "user-specific"
$ws.Range("D4").Value = "user-specific"
$ws.Range("F4").Value = "milligram"
$ws.Range("G4").Value = "UO"
$ws.Range("H4").Value = "http://purl.obolibrary.org/obo/UO_0000022"
$ws.Range("P4").Value = "microliter"
$ws.Range("Q4").Value = "UO"
$ws.Range("R4").Value = "http://purl.obolibrary.org/obo/UO_0000101"

# --- Row 5 (Protein -> protein, now with a real ontology term) ---
$ws.Range("B5").Value = "protein"
$ws.Range("C5").Value = "NFDI4PSO"
$ws.Range("D5").Value = "http://purl.obolibrary.org/obo/NFDI4PSO_1000093"
$ws.Range("F5").Value = "milligram"
$ws.Range("G5").Value = "UO"
$ws.Range("H5").Value = "http://purl.obolibrary.org/obo/UO_0000022"
$ws.Range("P5").Value = "microliter"
$ws.Range("Q5").Value = "UO"
$ws.Range("R5").Value = "http://purl.obolibrary.org/obo/UO_0000101"

# --- Bump version on metadata sheet ---
$meta = $wb.Worksheets.Item("SwateTemplateMetadata")
$meta.Range("B3").Value = "1.1.5"
